$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 (CON) values for columns B:E
$ws.Range("B2").Value = 59.371755917538131
$ws.Range("C2").Value = 46.75668948809188
$ws.Range("D2").Value = 59.438798947628541
$ws.Range("E2").Value = 49.966250868235662

# Update row 3 (STR) values for columns B:E
$ws.Range("B3").Value = 61.006834859898191
$ws.Range("C3").Value = 44.067825498757585
$ws.Range("D3").Value = 65.267752420178439
$ws.Range("E3").Value = 43.591477932075151

# Update the selection to match the new narrower range
$ws.Range("B1:E3").Select() | Out-Null
